# Update the cryptos list with latest prices / 1h volume changes
# (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to stay a text string even when the new value
    # looks like a number (e.g. "576.36"), matching the source data
    # which stores every Price/Volume value as plain text.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "62.798.96"
$ws.Range("E2").Value = "  +2.87%  "

Set-TextCell $ws.Range("D3") "2.444.31"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  -0.12%  "

Set-TextCell $ws.Range("D5") "576.36"
$ws.Range("E5").Value = "  +1.54%  "

Set-TextCell $ws.Range("D6") "145.75"
$ws.Range("E6").Value = "  +2.55%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("E8").Value = "  -0.15%  "

Set-TextCell $ws.Range("D9") "2.442.92"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  +2.52%  "

$ws.Range("E11").Value = "  +2.59%  "

$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("E14").Value = "  +6.99%  "

$ws.Range("E15").Value = "  +5.09%  "

Set-TextCell $ws.Range("D16") "2.888.65"
$ws.Range("E16").Value = "  +1.75%  "

Set-TextCell $ws.Range("D17") "62.672.76"
$ws.Range("E17").Value = "  +3.06%  "

Set-TextCell $ws.Range("D18") "2.449.73"
$ws.Range("E18").Value = "  +1.52%  "

Set-TextCell $ws.Range("D19") "7.92"
$ws.Range("E19").Value = "  -1.77%  "

Set-TextCell $ws.Range("D20") "11.04"
$ws.Range("E20").Value = "  +2.92%  "

Set-TextCell $ws.Range("D21") "330.41"
$ws.Range("E21").Value = "  +1.77%  "

$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("E23").Value = "  +7.27%  "

$ws.Range("E24").Value = "  +0.01%  "

Set-TextCell $ws.Range("D25") "66.37"
$ws.Range("E25").Value = "  +1.71%  "

Set-TextCell $ws.Range("D26") "647.69"
$ws.Range("E26").Value = "  +10.21%  "

$ws.Range("E27").Value = "  +17.96%  "

Set-TextCell $ws.Range("D28") "8.52"
$ws.Range("E28").Value = "  +2.89%  "

Set-TextCell $ws.Range("D29") "0.0$([char]0x2083)0990"
$ws.Range("E29").Value = "  +4.23%  "

Set-TextCell $ws.Range("D30") "2.566.37"
$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D31") "8.19"
$ws.Range("E31").Value = "  +1.76%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws.Range("D32") "1.44"
$ws.Range("E32").Value = "  +6.29%  "

$ws.Range("B33").Value = "BabyDogeCoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws.Range("D33") "0.0$([char]0x2086)0481"
$ws.Range("E33").Value = "  +67.86%  "

$ws.Range("E34").Value = "  +2.67%  "

$ws.Range("E35").Value = "  +3.99%  "

$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("E38").Value = "  +2.89%  "

$ws.Range("E39").Value = "  +5.73%  "

$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell $ws.Range("D40") "0.374"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D41") "153.33"
$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("E42").Value = "  +2.40%  "

$ws.Range("E43").Value = "  +7.54%  "

$ws.Range("E44").Value = "  +3.89%  "

Set-TextCell $ws.Range("D45") "42.53"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("E46").Value = "  +0.02%  "

Set-TextCell $ws.Range("D47") "14.95"
$ws.Range("E47").Value = "  +27.39%  "

Set-TextCell $ws.Range("D48") "145.21"
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("E49").Value = "  +2.93%  "

Set-TextCell $ws.Range("D50") "20.63"
$ws.Range("E50").Value = "  +5.28%  "

$ws.Range("E51").Value = "  +1.89%  "
